# Handles float input without breaking stuff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Summary rows 10-12: give the row-label cells the "mtitleStyle" look
# and refresh the score totals -------------------------------------------
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("B10").Value = 15
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = 10
$ws.Range("E10").Value = 28

$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("A12").Style = "mtitleStyle"
$ws.Range("B12").Value = 60
$ws.Range("C12").Value = -3
$ws.Range("E12").Value = "57/112"

# ---- Drop the third "Student Ans / Correct Ans" block (columns G:H) ----
$ws.Range("G15:H21").ClearContents()
$ws.Columns("G:H").Delete()

# ---- Drop the second block's per-question rows beyond 18 (D19:E40) -----
$ws.Range("D19:E40").ClearContents()

# ---- Fill in column A ("Student Ans") for the first block, grading it
# against column B ("Correct Ans") -----------------------------------
$ws.Range("A18").Value = "Option B"
$ws.Range("A18").Style = "correctStyle"
$ws.Range("A19").Value = "Option C"
$ws.Range("A19").Style = "correctStyle"
$ws.Range("A21").Value = "Option C"
$ws.Range("A21").Style = "correctStyle"
$ws.Range("A23").Value = "Option D"
$ws.Range("A23").Style = "correctStyle"
$ws.Range("A25").Value = "Option D"
$ws.Range("A25").Style = "incorrectStyle"
$ws.Range("A26").Value = "Option C"
$ws.Range("A26").Style = "correctStyle"
$ws.Range("A27").Value = "Option A"
$ws.Range("A27").Style = "correctStyle"
$ws.Range("A31").Value = "Option D"
$ws.Range("A31").Style = "correctStyle"
$ws.Range("A32").Value = "Option C"
$ws.Range("A32").Style = "correctStyle"
$ws.Range("A34").Value = "Option B"
$ws.Range("A34").Style = "correctStyle"
$ws.Range("A36").Value = "Option D"
$ws.Range("A36").Style = "incorrectStyle"
$ws.Range("A37").Value = "Option B"
$ws.Range("A37").Style = "incorrectStyle"
$ws.Range("A38").Value = "Option A"
$ws.Range("A38").Style = "correctStyle"
$ws.Range("A39").Value = "Option D"
$ws.Range("A39").Style = "correctStyle"
$ws.Range("A40").Value = "Option D"
$ws.Range("A40").Style = "correctStyle"

# ---- Column D ("Student Ans") for the second block, rows 16-18 only ----
$ws.Range("D16").Value = "Option A"
$ws.Range("D16").Style = "correctStyle"
$ws.Range("D17").Value = "Option C"
$ws.Range("D17").Style = "correctStyle"
$ws.Range("D18").Value = "Option D"
$ws.Range("D18").Style = "correctStyle"
